$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "R pour revenir..." instruction text in B2 to use "Espace" instead of "R"
$ws.Range("B2").Value = "Espace pour revenir au début des instructions`nBienvenue à la course des champions `nMerci d'avoir rejoint le jury !"

# Reset the view: scroll back to top-left and select B2 (instead of A3/B8)
$ws.Range("A1").Select()
$ws.Range("B2").Select()
